$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old data range (A1:G2) before laying out the new set.
$ws.Range("A1:H3").Clear()

# Header row
$ws.Range("A1").Value = "x1"
$ws.Range("B1").Value = "y1"
$ws.Range("C1").Value = "stim1_color"
$ws.Range("D1").Value = "probe1_color"
$ws.Range("E1").Value = "x2"
$ws.Range("F1").Value = "y2"
$ws.Range("G1").Value = "stim2_color"
$ws.Range("H1").Value = "probe2_color"

# Row 2
$ws.Range("A2").Value = -0.25
$ws.Range("B2").Value = -0.25
$ws.Range("C2").Value = "red"
$ws.Range("D2").Value = "blue"
$ws.Range("E2").Value = -0.25
$ws.Range("F2").Value = -0.25
$ws.Range("G2").Value = "orange"
$ws.Range("H2").Value = "orange"

# Row 3
$ws.Range("A3").Value = -0.25
$ws.Range("B3").Value = -0.25
$ws.Range("C3").Value = "orange"
$ws.Range("D3").Value = "orange"
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 0.25
$ws.Range("G3").Value = "blue"
$ws.Range("H3").Value = "green"

# Update selection to match target state
$ws.Range("F7").Select()
